$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.846.54'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '1.638.05'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').Value = '''0.5060'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').Value = '''0.2577'
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('D9').Value = '''0.06426'
$ws.Range('E9').Value = '  +1.07%  '
$ws.Range('E10').Value = '  -0.70%  '
$ws.Range('D11').Value = '''0.07767'
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').Value = '''4.282'
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('D13').Value = '1.864.22'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').Value = '1.635.32'
$ws.Range('E14').Value = '  -0.10%  '
$ws.Range('D15').Value = '''0.5635'
$ws.Range('E15').Value = '  +3.66%  '
$ws.Range('D16').Value = '0.0₅7597'
$ws.Range('E16').Value = '  -1.74%  '
$ws.Range('D17').Value = '''63.11'
$ws.Range('E17').Value = '  -1.47%  '
$ws.Range('D18').Value = '25.862.76'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').Value = '''194.74'
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').Value = '''4.318'
$ws.Range('E21').Value = '  -2.56%  '
$ws.Range('D22').Value = '''9.873'
$ws.Range('E22').Value = '  -0.38%  '
$ws.Range('D23').Value = '''6.093'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '''1.798'
$ws.Range('E25').Value = '  -4.84%  '
$ws.Range('D26').Value = '''0.1272'
$ws.Range('E26').Value = '  +2.26%  '
$ws.Range('D27').Value = '''139.86'
$ws.Range('E27').Value = '  -2.19%  '
$ws.Range('D28').Value = '''6.793'
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('D29').Value = '''15.45'
$ws.Range('E29').Value = '  -1.00%  '
$ws.Range('D30').Value = '''1.243'
$ws.Range('E30').Value = '  +0.58%  '
$ws.Range('D31').Value = '''0.04866'
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('D32').Value = '''3.298'
$ws.Range('E32').Value = '  +1.81%  '
$ws.Range('D33').Value = '''3.218'
$ws.Range('E33').Value = '  +0.72%  '
$ws.Range('D34').Value = '''1.559'
$ws.Range('E34').Value = '  +0.57%  '
$ws.Range('D35').Value = '''2.377'
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('D36').Value = '''0.9036'
$ws.Range('E36').Value = '  -0.79%  '
$ws.Range('D38').Value = '1.130.51'
$ws.Range('E38').Value = '  +0.45%  '
$ws.Range('D39').Value = '''0.5513'
$ws.Range('E39').Value = '  +0.34%  '
$ws.Range('D40').Value = '''0.01560'
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('D41').Value = '''0.9962'
$ws.Range('E41').Value = '  -0.47%  '
$ws.Range('D42').Value = '''5.524'
$ws.Range('E42').Value = '  -0.97%  '
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('D44').Value = '''97.84'
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('D45').Value = '1.774.60'
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('D46').Value = '0.0₈113'
$ws.Range('E46').Value = '  -8.56%  '
$ws.Range('D47').Value = '''55.37'
$ws.Range('E47').Value = '  +0.65%  '
$ws.Range('D48').Value = '''0.4368'
$ws.Range('E48').Value = '  -2.54%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '''7.705'
$ws.Range('E49').Value = '  +2.43%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.05048'
$ws.Range('E50').Value = '  -2.14%  '
$ws.Range('D51').Value = '''1.004'
$ws.Range('E51').Value = '  +0.16%  '
